$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = "OR.0008.0031.20230907.16"
$ws.Range("A10").Value = "OR.0008.0031.20230907.17"
